{"js": "// Replace every arithmetic expression in the worksheet's 20x5 table with the\n// new expression from the commit (each of the 100 cells changes, in order).\nconst newValues = [\n  [\"42-32=\", \"6+62=\", \"20+40=\", \"25-16=\", \"25+10=\"],\n  [\"12-1=\", \"68-51=\", \"63+20=\", \"38-33=\", \"96-39=\"],\n  [\"36-9=\", \"20+56=\", \"9+58=\", \"86-65=\", \"24+54=\"],\n  [\"64+20=\", \"8+36=\", \"84-27=\", \"94-10=\", \"68-45=\"],\n  [\"60-53=\", \"62+1=\", \"7+85=\", \"71-66=\", \"53+0=\"],\n  [\"72-32=\", \"92-60=\", \"44-17=\", \"32+0=\", \"26-16=\"],\n  [\"22+24=\", \"80-60=\", \"44+20=\", \"32+22=\", \"49-19=\"],\n  [\"54+44=\", \"66+19=\", \"28+37=\", \"17-17=\", \"89+7=\"],\n  [\"70-21=\", \"67-40=\", \"11+87=\", \"60+0=\", \"32+17=\"],\n  [\"59+5=\", \"87-8=\", \"5+30=\", \"36+47=\", \"82-40=\"],\n  [\"11+75=\", \"75+22=\", \"43-13=\", \"1+41=\", \"29+11=\"],\n  [\"81-28=\", \"78-7=\", \"63-42=\", \"62-45=\", \"54-37=\"],\n  [\"71+7=\", \"64+8=\", \"4+43=\", \"82-56=\", \"43-34=\"],\n  [\"74+15=\", \"84+14=\", \"16+28=\", \"3+3=\", \"80+4=\"],\n  [\"56-6=\", \"51+31=\", \"83-54=\", \"33-9=\", \"86-52=\"],\n  [\"51+27=\", \"81-30=\", \"38+48=\", \"87-77=\", \"8+37=\"],\n  [\"97-61=\", \"25-22=\", \"39+28=\", \"70-10=\", \"96+0=\"],\n  [\"52+13=\", \"93-87=\", \"15+34=\", \"39-20=\", \"78-2=\"],\n  [\"19+25=\", \"74-49=\", \"84-68=\", \"55-38=\", \"22+6=\"],\n  [\"83-67=\", \"77-72=\", \"62+7=\", \"38+34=\", \"2+26=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace every arithmetic expression in the worksheet's 20x5 table with the\n# new expression from the commit (each of the 100 cells changes, in order).\n$newValues = @(\n    @(\"42-32=\", \"6+62=\", \"20+40=\", \"25-16=\", \"25+10=\"),\n    @(\"12-1=\", \"68-51=\", \"63+20=\", \"38-33=\", \"96-39=\"),\n    @(\"36-9=\", \"20+56=\", \"9+58=\", \"86-65=\", \"24+54=\"),\n    @(\"64+20=\", \"8+36=\", \"84-27=\", \"94-10=\", \"68-45=\"),\n    @(\"60-53=\", \"62+1=\", \"7+85=\", \"71-66=\", \"53+0=\"),\n    @(\"72-32=\", \"92-60=\", \"44-17=\", \"32+0=\", \"26-16=\"),\n    @(\"22+24=\", \"80-60=\", \"44+20=\", \"32+22=\", \"49-19=\"),\n    @(\"54+44=\", \"66+19=\", \"28+37=\", \"17-17=\", \"89+7=\"),\n    @(\"70-21=\", \"67-40=\", \"11+87=\", \"60+0=\", \"32+17=\"),\n    @(\"59+5=\", \"87-8=\", \"5+30=\", \"36+47=\", \"82-40=\"),\n    @(\"11+75=\", \"75+22=\", \"43-13=\", \"1+41=\", \"29+11=\"),\n    @(\"81-28=\", \"78-7=\", \"63-42=\", \"62-45=\", \"54-37=\"),\n    @(\"71+7=\", \"64+8=\", \"4+43=\", \"82-56=\", \"43-34=\"),\n    @(\"74+15=\", \"84+14=\", \"16+28=\", \"3+3=\", \"80+4=\"),\n    @(\"56-6=\", \"51+31=\", \"83-54=\", \"33-9=\", \"86-52=\"),\n    @(\"51+27=\", \"81-30=\", \"38+48=\", \"87-77=\", \"8+37=\"),\n    @(\"97-61=\", \"25-22=\", \"39+28=\", \"70-10=\", \"96+0=\"),\n    @(\"52+13=\", \"93-87=\", \"15+34=\", \"39-20=\", \"78-2=\"),\n    @(\"19+25=\", \"74-49=\", \"84-68=\", \"55-38=\", \"22+6=\"),\n    @(\"83-67=\", \"77-72=\", \"62+7=\", \"38+34=\", \"2+26=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
